# Adds 15 duplicated match rows (17-31) to the "Manish Pandey" sheet,
# mirroring rows 2-16 in a different order, and extends the used range
# + ignoredErrors sqref from A1:K16 to A1:K31.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new rows to be stored as text (matching the existing sheet,
# where numeric-looking values like scores/strike-rates are text, not numbers).
$ws.Range("A17:K31").NumberFormat = "@"

# Row 17
$ws.Range("A17").Value = ' Dubai (DSC)'
$ws.Range("B17").Value = ' October 13 2020'
$ws.Range("C17").Value = 'Super Kings won by 20 runs'
$ws.Range("D17").Value = 'Sunrisers Hyderabad'
$ws.Range("E17").Value = 'Chennai Super Kings'
$ws.Range("F17").Value = 'Manish Pandey '
$ws.Range("G17").Value = '4'
$ws.Range("H17").Value = '3'
$ws.Range("I17").Value = '1'
$ws.Range("J17").Value = '0'
$ws.Range("K17").Value = '133.33'

# Row 18
$ws.Range("A18").Value = ' Dubai (DSC)'
$ws.Range("B18").Value = ' October 22 2020'
$ws.Range("C18").Value = 'Sunrisers won by 8 wickets (with 11 balls remaining)'
$ws.Range("D18").Value = 'Sunrisers Hyderabad'
$ws.Range("E18").Value = 'Rajasthan Royals'
$ws.Range("F18").Value = 'Manish Pandey '
$ws.Range("G18").Value = '83'
$ws.Range("H18").Value = '47'
$ws.Range("I18").Value = '4'
$ws.Range("J18").Value = '8'
$ws.Range("K18").Value = '176.59'

# Row 19
$ws.Range("A19").Value = ' Abu Dhabi'
$ws.Range("B19").Value = ' October 18 2020'
$ws.Range("C19").Value = 'Match tied (KKR won the one-over eliminator)'
$ws.Range("D19").Value = 'Sunrisers Hyderabad'
$ws.Range("E19").Value = 'Kolkata Knight Riders'
$ws.Range("F19").Value = 'Manish Pandey '
$ws.Range("G19").Value = '6'
$ws.Range("H19").Value = '7'
$ws.Range("I19").Value = '0'
$ws.Range("J19").Value = '0'
$ws.Range("K19").Value = '85.71'

# Row 20
$ws.Range("A20").Value = ' Dubai (DSC)'
$ws.Range("B20").Value = ' October 27 2020'
$ws.Range("C20").Value = 'Sunrisers won by 88 runs'
$ws.Range("D20").Value = 'Sunrisers Hyderabad'
$ws.Range("E20").Value = 'Delhi Capitals'
$ws.Range("F20").Value = 'Manish Pandey '
$ws.Range("G20").Value = '44'
$ws.Range("H20").Value = '31'
$ws.Range("I20").Value = '4'
$ws.Range("J20").Value = '1'
$ws.Range("K20").Value = '141.93'

# Row 21
$ws.Range("A21").Value = ' Dubai (DSC)'
$ws.Range("B21").Value = ' October 02 2020'
$ws.Range("C21").Value = 'Sunrisers won by 7 runs'
$ws.Range("D21").Value = 'Sunrisers Hyderabad'
$ws.Range("E21").Value = 'Chennai Super Kings'
$ws.Range("F21").Value = 'Manish Pandey '
$ws.Range("G21").Value = '29'
$ws.Range("H21").Value = '21'
$ws.Range("I21").Value = '5'
$ws.Range("J21").Value = '0'
$ws.Range("K21").Value = '138.09'

# Row 22
$ws.Range("A22").Value = ' Abu Dhabi'
$ws.Range("B22").Value = ' September 29 2020'
$ws.Range("C22").Value = 'Sunrisers won by 15 runs'
$ws.Range("D22").Value = 'Sunrisers Hyderabad'
$ws.Range("E22").Value = 'Delhi Capitals'
$ws.Range("F22").Value = 'Manish Pandey '
$ws.Range("G22").Value = '3'
$ws.Range("H22").Value = '5'
$ws.Range("I22").Value = '0'
$ws.Range("J22").Value = '0'
$ws.Range("K22").Value = '60.00'

# Row 23
$ws.Range("A23").Value = ' Sharjah'
$ws.Range("B23").Value = ' October 04 2020'
$ws.Range("C23").Value = 'Mumbai won by 34 runs'
$ws.Range("D23").Value = 'Sunrisers Hyderabad'
$ws.Range("E23").Value = 'Mumbai Indians'
$ws.Range("F23").Value = 'Manish Pandey '
$ws.Range("G23").Value = '30'
$ws.Range("H23").Value = '19'
$ws.Range("I23").Value = '4'
$ws.Range("J23").Value = '1'
$ws.Range("K23").Value = '157.89'

# Row 24
$ws.Range("A24").Value = ' Abu Dhabi'
$ws.Range("B24").Value = ' September 26 2020'
$ws.Range("C24").Value = 'KKR won by 7 wickets (with 12 balls remaining)'
$ws.Range("D24").Value = 'Sunrisers Hyderabad'
$ws.Range("E24").Value = 'Kolkata Knight Riders'
$ws.Range("F24").Value = 'Manish Pandey '
$ws.Range("G24").Value = '51'
$ws.Range("H24").Value = '38'
$ws.Range("I24").Value = '3'
$ws.Range("J24").Value = '2'
$ws.Range("K24").Value = '134.21'

# Row 25
$ws.Range("A25").Value = ' Abu Dhabi'
$ws.Range("B25").Value = ' November 06 2020'
$ws.Range("C25").Value = 'Sunrisers won by 6 wickets (with 2 balls remaining)'
$ws.Range("D25").Value = 'Sunrisers Hyderabad'
$ws.Range("E25").Value = 'Royal Challengers Bangalore'
$ws.Range("F25").Value = 'Manish Pandey '
$ws.Range("G25").Value = '24'
$ws.Range("H25").Value = '21'
$ws.Range("I25").Value = '3'
$ws.Range("J25").Value = '1'
$ws.Range("K25").Value = '114.28'

# Row 26
$ws.Range("A26").Value = ' Sharjah'
$ws.Range("B26").Value = ' October 31 2020'
$ws.Range("C26").Value = 'Sunrisers won by 5 wickets (with 35 balls remaining)'
$ws.Range("D26").Value = 'Sunrisers Hyderabad'
$ws.Range("E26").Value = 'Royal Challengers Bangalore'
$ws.Range("F26").Value = 'Manish Pandey '
$ws.Range("G26").Value = '26'
$ws.Range("H26").Value = '19'
$ws.Range("I26").Value = '3'
$ws.Range("J26").Value = '1'
$ws.Range("K26").Value = '136.84'

# Row 27
$ws.Range("A27").Value = ' Dubai (DSC)'
$ws.Range("B27").Value = ' September 21 2020'
$ws.Range("C27").Value = 'RCB won by 10 runs'
$ws.Range("D27").Value = 'Sunrisers Hyderabad'
$ws.Range("E27").Value = 'Royal Challengers Bangalore'
$ws.Range("F27").Value = 'Manish Pandey '
$ws.Range("G27").Value = '34'
$ws.Range("H27").Value = '33'
$ws.Range("I27").Value = '3'
$ws.Range("J27").Value = '1'
$ws.Range("K27").Value = '103.03'

# Row 28
$ws.Range("A28").Value = ' Abu Dhabi'
$ws.Range("B28").Value = ' November 08 2020'
$ws.Range("C28").Value = 'Capitals won by 17 runs'
$ws.Range("D28").Value = 'Sunrisers Hyderabad'
$ws.Range("E28").Value = 'Delhi Capitals'
$ws.Range("F28").Value = 'Manish Pandey '
$ws.Range("G28").Value = '21'
$ws.Range("H28").Value = '14'
$ws.Range("I28").Value = '3'
$ws.Range("J28").Value = '0'
$ws.Range("K28").Value = '150.00'

# Row 29
$ws.Range("A29").Value = ' Dubai (DSC)'
$ws.Range("B29").Value = ' October 11 2020'
$ws.Range("C29").Value = 'Royals won by 5 wickets (with 1 ball remaining)'
$ws.Range("D29").Value = 'Sunrisers Hyderabad'
$ws.Range("E29").Value = 'Rajasthan Royals'
$ws.Range("F29").Value = 'Manish Pandey '
$ws.Range("G29").Value = '54'
$ws.Range("H29").Value = '44'
$ws.Range("I29").Value = '2'
$ws.Range("J29").Value = '3'
$ws.Range("K29").Value = '122.72'

# Row 30
$ws.Range("A30").Value = ' Dubai (DSC)'
$ws.Range("B30").Value = ' October 24 2020'
$ws.Range("C30").Value = 'Kings XI won by 12 runs'
$ws.Range("D30").Value = 'Sunrisers Hyderabad'
$ws.Range("E30").Value = 'Kings XI Punjab'
$ws.Range("F30").Value = 'Manish Pandey '
$ws.Range("G30").Value = '15'
$ws.Range("H30").Value = '29'
$ws.Range("I30").Value = '0'
$ws.Range("J30").Value = '0'
$ws.Range("K30").Value = '51.72'

# Row 31
$ws.Range("A31").Value = ' Dubai (DSC)'
$ws.Range("B31").Value = ' October 08 2020'
$ws.Range("C31").Value = 'Sunrisers won by 69 runs'
$ws.Range("D31").Value = 'Sunrisers Hyderabad'
$ws.Range("E31").Value = 'Kings XI Punjab'
$ws.Range("F31").Value = 'Manish Pandey '
$ws.Range("G31").Value = '1'
$ws.Range("H31").Value = '2'
$ws.Range("I31").Value = '0'
$ws.Range("J31").Value = '0'
$ws.Range("K31").Value = '50.00'
